$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The statement gained one more overdue period (2509). Insert a new row
# right after the last detail row (59) and before the trailing blank
# rows that lead into the signature block. This pushes the old rows
# 64/65 (signature lines) down to 65/66.
$ws.Rows.Item(60).Insert()

# Row 59 used to be the last row of the detail table and therefore
# carried the "closing" bottom-border styling. Move that styling down
# onto the freshly inserted row 60 (the new last row), and restore the
# regular interior-row styling (copied from row 58) onto row 59.
$ws.Range("B59:J59").Copy()
$ws.Range("B60:J60").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B58:J58").Copy()
$ws.Range("B59:J59").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new detail row for period 2509.
$ws.Range("B60").Value = "CC"
$ws.Range("C60").Value = "1047424235"
$ws.Range("D60").Value = "YISETH MARGARITA MORALES BOGALLO"
$ws.Range("E60").Value = "2509"
$ws.Range("F60").Value = 36341
$ws.Range("G60").Value = 908526

# Roll the new period into the account totals.
$ws.Range("E11").Value = 1635345
$ws.Range("F13").Value = 45
